$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the column-A "Model Level" style down onto the new row 5 before writing values
$ws.Range("A4").Copy() | Out-Null
$ws.Range("A5").PasteSpecial(-4122) | Out-Null

# Clear every text cell in the results table (even ones whose value is unchanged) so the
# shared-string table is rebuilt from empty, in exactly the left-to-right, top-to-bottom
# order the values are written below -- this keeps engine-assigned string ids deterministic.
$ws.Range("C2:T5").ClearContents()

# Row 2
$ws.Range("A2").Value = [double]"0"
$ws.Range("B2").Value = [double]"1"
$ws.Range("C2").Value = '[''N1ratio-ArgsPreds'']'
$ws.Range("D2").Value = [double]"1546"
$ws.Range("E2").Value = [double]"1544"
$ws.Range("F2").Value = [double]"1"
$ws.Range("G2").Value = [double]"0.2788887100714704"
$ws.Range("H2").Value = [double]"597.139684767698"
$ws.Range("I2").Value = [double]"9.159653958517076e-112"
$ws.Range("J2").Value = [double]"172.1361681705525"
$ws.Range("K2").Value = [double]"238.7095730918499"
$ws.Range("L2").Value = [double]"66.5734049212974"
$ws.Range("M2").Value = [double]"0.111487155550876"
$ws.Range("N2").Value = [double]"0.1545045780529773"
$ws.Range("O2").Value = '{''const'': 0.7216232857673603, ''N1ratio-ArgsPreds'': -0.20623521687888982}'
$ws.Range("P2").Value = '{''const'': 5.823844912460471e-164, ''N1ratio-ArgsPreds'': 9.159653958519536e-112}'
$ws.Range("Q2").Value = '{''N1ratio-ArgsPreds'': -0.5280991479556377}'
$ws.Range("R2").Value = '{''N1ratio-ArgsPreds'': -0.5280991479556364}'
$ws.Range("S2").Value = '{''N1ratio-ArgsPreds'': -0.5280991479556364}'
$ws.Range("T2").Value = '{''N1ratio-ArgsPreds'': 27.888871007146914}'

# Row 3
$ws.Range("A3").Value = [double]"1"
$ws.Range("B3").Value = [double]"2"
$ws.Range("C3").Value = '[''N1ratio-ArgsPreds'', ''latitude'', ''longitude'', ''Macro_class'']'
$ws.Range("D3").Value = [double]"1546"
$ws.Range("E3").Value = [double]"1541"
$ws.Range("F3").Value = [double]"4"
$ws.Range("G3").Value = [double]"0.3351199849338026"
$ws.Range("H3").Value = [double]"194.1778535528601"
$ws.Range("I3").Value = [double]"6.876195971999602e-135"
$ws.Range("J3").Value = [double]"158.7132245537547"
$ws.Range("K3").Value = [double]"238.7095730918499"
$ws.Range("L3").Value = [double]"19.9990871345238"
$ws.Range("M3").Value = [double]"0.1029936564268363"
$ws.Range("N3").Value = [double]"0.1545045780529773"
$ws.Range("O3").Value = '{''const'': 0.5920993353367638, ''N1ratio-ArgsPreds'': -0.21008109385060292, ''latitude'': 0.003788815465326568, ''longitude'': -4.566087073707778e-06, ''Macro_class'': 0.0451128018370466}'
$ws.Range("P3").Value = '{''const'': 1.1460070432133304e-100, ''N1ratio-ArgsPreds'': 2.7241569733103175e-112, ''latitude'': 2.010532168931972e-14, ''longitude'': 0.9648727161654734, ''Macro_class'': 1.810321629638847e-21}'
$ws.Range("Q3").Value = '{''N1ratio-ArgsPreds'': -0.5379471476457032, ''latitude'': 0.16892586870392384, ''longitude'': -0.0009843143319924143, ''Macro_class'': 0.2056838925374642}'
$ws.Range("R3").Value = '{''N1ratio-ArgsPreds'': -0.5296028099294897, ''latitude'': 0.193070598681436, ''longitude'': -0.0011220543053802913, ''Macro_class'': 0.23891601510265767}'
$ws.Range("S3").Value = '{''N1ratio-ArgsPreds'': -0.5090966326675047, ''latitude'': 0.16044897288356422, ''longitude'': -0.0009149256209509112, ''Macro_class'': 0.20062251798959566}'
$ws.Range("T3").Value = '{''N1ratio-ArgsPreds'': 25.917938139339224, ''latitude'': 2.5743872899390725, ''longitude'': 8.370888918724104e-05, ''Macro_class'': 4.024939472448564}'
$ws.Range("U3").Value = [double]"0.05623127486233215"
$ws.Range("V3").Value = [double]"43.44262252802935"
$ws.Range("W3").Value = [double]"6.001493159104878e-27"

# Row 4
$ws.Range("A4").Value = [double]"2"
$ws.Range("B4").Value = [double]"3"
$ws.Range("C4").Value = '[''N1ratio-ArgsPreds'', ''latitude'', ''longitude'', ''Macro_class'', ''Fam_class'']'
$ws.Range("D4").Value = [double]"1546"
$ws.Range("E4").Value = [double]"1540"
$ws.Range("F4").Value = [double]"5"
$ws.Range("G4").Value = [double]"0.3358954086663867"
$ws.Range("H4").Value = [double]"155.7823680476198"
$ws.Range("I4").Value = [double]"4.170842569285563e-134"
$ws.Range("J4").Value = [double]"158.5281234855843"
$ws.Range("K4").Value = [double]"238.7095730918499"
$ws.Range("L4").Value = [double]"16.03628992125313"
$ws.Range("M4").Value = [double]"0.1029403399257041"
$ws.Range("N4").Value = [double]"0.1545045780529773"
$ws.Range("O4").Value = '{''const'': 0.5986196932832282, ''N1ratio-ArgsPreds'': -0.20686377387951788, ''latitude'': 0.0038003308375348654, ''longitude'': -4.031102623384531e-05, ''Macro_class'': 0.046025152341286533, ''Fam_class'': -0.00032243452209890574}'
$ws.Range("P4").Value = '{''const'': 1.3719407742700062e-99, ''N1ratio-ArgsPreds'': 1.0482111701943141e-102, ''latitude'': 1.6734942419767786e-14, ''longitude'': 0.7064475226020194, ''Macro_class'': 7.538214510784381e-22, ''Fam_class'': 0.18013535613409756}'
$ws.Range("Q4").Value = '{''N1ratio-ArgsPreds'': -0.5297086713992897, ''latitude'': 0.16943928622756352, ''longitude'': -0.0086898738939457, ''Macro_class'': 0.2098435943389268, ''Fam_class'': -0.029619705700312356}'
$ws.Range("R4").Value = '{''N1ratio-ArgsPreds'': -0.5096588369434795, ''latitude'': 0.19371561511617877, ''longitude'': -0.009598794712681876, ''Macro_class'': 0.24120124994063472, ''Fam_class'': -0.034150565916667}'
$ws.Range("S4").Value = '{''N1ratio-ArgsPreds'': -0.4827354770349507, ''latitude'': 0.16091195426469926, ''longitude'': -0.007822668978782919, ''Macro_class'': 0.20254120521153235, ''Fam_class'': -0.027846431236075636}'
$ws.Range("T4").Value = '{''N1ratio-ArgsPreds'': 23.303354078816138, ''latitude'': 2.5892657025284667, ''longitude'': 0.006119414995161259, ''Macro_class'': 4.1022939808540055, ''Fam_class'': 0.07754237325854889}'
$ws.Range("U4").Value = [double]"0.0007754237325841018"
$ws.Range("V4").Value = [double]"1.798139274690895"
$ws.Range("W4").Value = [double]"0.1801353561342204"

# Row 5
$ws.Range("A5").Value = [double]"3"
$ws.Range("B5").Value = [double]"4"
$ws.Range("C5").Value = '[''N1ratio-ArgsPreds'', ''latitude'', ''longitude'', ''Macro_class'', ''Fam_class'', ''Nlen_freq'', ''Vlen_freq'']'
$ws.Range("D5").Value = [double]"1546"
$ws.Range("E5").Value = [double]"1538"
$ws.Range("F5").Value = [double]"7"
$ws.Range("G5").Value = [double]"0.3450160878646346"
$ws.Range("H5").Value = [double]"115.7356110594802"
$ws.Range("I5").Value = [double]"1.676257582139428e-136"
$ws.Range("J5").Value = [double]"156.3509300478628"
$ws.Range("K5").Value = [double]"238.7095730918499"
$ws.Range("L5").Value = [double]"11.7655204348553"
$ws.Range("M5").Value = [double]"0.1016586021117444"
$ws.Range("N5").Value = [double]"0.1545045780529773"
$ws.Range("O5").Value = '{''const'': 0.6613995821693961, ''N1ratio-ArgsPreds'': -0.20022658501321183, ''latitude'': 0.003904782868167093, ''longitude'': -0.00013919989562284846, ''Macro_class'': 0.041249893628343776, ''Fam_class'': -0.0005341528657763941, ''Nlen_freq'': -0.05245277165013973, ''Vlen_freq'': 0.04521620320763898}'
$ws.Range("P5").Value = '{''const'': 3.4679513017458327e-38, ''N1ratio-ArgsPreds'': 2.6913691577768e-95, ''latitude'': 2.5250829575406025e-15, ''longitude'': 0.20243259343697534, ''Macro_class'': 3.541864495018582e-16, ''Fam_class'': 0.029582124451356717, ''Nlen_freq'': 1.099034877332113e-05, ''Vlen_freq'': 1.0268744561068315e-05}'
$ws.Range("Q5").Value = '{''N1ratio-ArgsPreds'': -0.5127130591165668, ''latitude'': 0.1740963222257306, ''longitude'': -0.030007411173207665, ''Macro_class'': 0.18807164136869625, ''Fam_class'': -0.04906872434218462, ''Nlen_freq'': -0.16839277217585, ''Vlen_freq'': 0.17778809172835128}'
$ws.Range("R5").Value = '{''N1ratio-ArgsPreds'': -0.4934305412228254, ''latitude'': 0.19975306738554952, ''longitude'': -0.032498939063819676, ''Macro_class'': 0.20570605550116072, ''Fam_class'': -0.055442836886884554, ''Nlen_freq'': -0.1117780702482574, ''Vlen_freq'': 0.11215005622049709}'
$ws.Range("S5").Value = '{''N1ratio-ArgsPreds'': -0.4591232942458281, ''latitude'': 0.1649874194744548, ''longitude'': -0.026315641903484864, ''Macro_class'': 0.17011830467528202, ''Fam_class'': -0.04493960793348434, ''Nlen_freq'': -0.09103368573944236, ''Vlen_freq'': 0.09134048898367954}'
$ws.Range("T5").Value = '{''N1ratio-ArgsPreds'': 21.079419931914124, ''latitude'': 2.722084858483971, ''longitude'': 0.06925130087924485, ''Macro_class'': 2.894023758559208, ''Fam_class'': 0.20195683612152884, ''Nlen_freq'': 0.8287131939307552, ''Vlen_freq'': 0.8343084927777683}'
$ws.Range("U5").Value = [double]"0.009120679198247883"
$ws.Range("V5").Value = [double]"10.7083581344561"
$ws.Range("W5").Value = [double]"2.407140461331682e-05"
